$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.701.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.92%  "

$ws.Range("E9").Value = "  +1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.863.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.674.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.38%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.96%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.122"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.98%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("E35").Value = "  -1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.198.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("E37").Value = "  +5.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.506"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.796"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.772.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("E46").Value = "  +2.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.411"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("E51").Value = "  -0.07%  "
